$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D: re-affirm age formula (keeps same text/values) ---
$ws.Range("D3").Formula = '=INT((TODAY()-C3)/365)'

# --- Column E: "Юбилей" flag ---
$ws.Range("E3").Formula = '=IF(MOD(D3,5)=0,"Юбилей","")'
$ws.Range("E4:E12").Formula = '=IF(MOD(D4,5)=0,"Юбилей","")'

# --- Column F: bonus amount, currency-style formatted ---
$ws.Range("F3").Formula = '=IF(E3="Юбилей",50,0)'
$ws.Range("F4:F14").Formula = '=IF(E4="Юбилей",50,0)'

$fmt = '_-[$$-409]* #,##0.00_ ;_-[$$-409]* \-#,##0.00\ ;_-[$$-409]* "-"??_ ;_-@_ '
$ws.Range("F3:F14").NumberFormat = $fmt

# --- Row 13: count of jubilees ---
$ws.Range("B13").Value = "Количествое юбиляров"
$ws.Range("B13:D13").Merge()
$ws.Range("E13").Formula = '=COUNTIF(E3:E12,"Юбилей")'

# --- Row 14: count of people younger than 30 ---
$ws.Range("B14").Value = "Количество человек моложе 30 лет"
$ws.Range("B14:D14").Merge()
$ws.Range("E14").Formula = '=COUNTIF(D3:D12, "< 30")'

# center-align B13:D14 (matches style used for header row / blank cells)
$ws.Range("B13:D14").HorizontalAlignment = -4108

# clear the leftover F13/F14 content but keep the currency style applied above
$ws.Range("F13:F14").ClearContents()

# --- conditional formatting: bold+italic when E = "Юбилей" ---
$rng = $ws.Range("E3:E14")
$fc = $rng.FormatConditions.Add(1, 3, '"Юбилей"')
$fc.Font.Bold = $true
$fc.Font.Italic = $true
$fc.Font.ThemeColor = 3

# --- column F width ---
$ws.Columns("F").ColumnWidth = 8.716666666666667

# --- move the active selection like the saved workbook shows ---
$ws.Range("F30").Select()
